$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.633.80"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.29%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.091.37"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.99%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.008"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.53%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "342.31"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.49%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.007"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.38%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5161"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.88%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4398"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.44%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09249"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.55%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "51.94"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.44%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.178"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.57%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "25.10"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.46%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.085.29"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.91%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.744"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.20%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.189"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.99%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "100.18"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.48%  "

$ws.Range("E17").Value = "  -1.57%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.008"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.52%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "21.14"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +9.19%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.06635"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.08%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.007"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.36%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.185"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.69%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "29.685.01"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.34%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.66"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.37%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.311"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.98%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.335.01"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.74%  "

$ws.Range("E27").Value = "  -2.04%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "162.85"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.65%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.522"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.24%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "132.61"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.79%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.135"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.39%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1051"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.17%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.631"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.23%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.183"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.66%  "

$ws.Range("E35").Value = "  -1.36%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.032"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.34%  "

$ws.Range("E37").Value = "  +1.22%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02568"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.82%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06703"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.78%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "12.48"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.64%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2244"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.82%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6828"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.70%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.294"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.87%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6614"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.00%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.06"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.88%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.307"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.61%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.605"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.28%  "

$ws.Range("E49").Value = "  -5.08%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "81.56"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.21%  "

$ws.Range("E51").Value = "  -2.11%  "
